# Scheduled market-price refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) for the affected leve rows across all 8 crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Beeswax
$ws.Range("H12").Value = 398
$ws.Range("I12").Value = 398
$ws.Range("K12").Value = 398
$ws.Range("M12").Value = -228

# Row 39: Hi-Potion of Mind
$ws.Range("H39").Value = 2874.625
$ws.Range("I39").Value = 2331.5
$ws.Range("K39").Value = 6994.5
$ws.Range("M39").Value = -6698.5

# Row 88: Growth Formula Zeta
$ws.Range("H88").Value = 1313.1666
$ws.Range("I88").Value = 1824.6
$ws.Range("J88").Value = 1116.4615
$ws.Range("K88").Value = 1824.6
$ws.Range("L88").Value = 1116.4615
$ws.Range("M88").Value = -1418.6
$ws.Range("N88").Value = -1928.4615

# Row 91: Growth Formula Zeta
$ws.Range("H91").Value = 1313.1666
$ws.Range("I91").Value = 1824.6
$ws.Range("J91").Value = 1116.4615
$ws.Range("K91").Value = 1824.6
$ws.Range("L91").Value = 1116.4615
$ws.Range("M91").Value = -420.5999999999999
$ws.Range("N91").Value = -3924.4615

# Row 98: Enchanted Durium Ink
$ws.Range("H98").Value = 1330
$ws.Range("I98").Value = 1412.5
$ws.Range("K98").Value = 1412.5
$ws.Range("M98").Value = 85.5

# Row 113: Starch Glue
$ws.Range("H113").Value = 4200
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -1246

# Row 116: Growth Formula Kappa
$ws.Range("H116").Value = 4865.778
$ws.Range("I116").Value = 3849
$ws.Range("K116").Value = 3849
$ws.Range("M116").Value = -407

# Row 122: Enchanted High Durium Ink
$ws.Range("H122").Value = 1330
$ws.Range("I122").Value = 1412.5
$ws.Range("K122").Value = 4237.5
$ws.Range("M122").Value = -1787.5

$ws = $wb.Worksheets.Item("ARM")
# Row 3: Bronze Skillet
$ws.Range("H3").Value = 2500
$ws.Range("J3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("N3").Value = -2730

# Row 17: Amateur's Skillet
$ws.Range("H17").Value = 17499.75
$ws.Range("I17").Value = 15000
$ws.Range("J17").Value = 19999.5
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 19999.5
$ws.Range("M17").Value = -14827
$ws.Range("N17").Value = -20345.5

# Row 32: Steel Ingot
$ws.Range("H32").Value = 4533354.5
$ws.Range("I32").Value = 4379189
$ws.Range("K32").Value = 4379189
$ws.Range("M32").Value = -4378902

# Row 36: Heavy Iron Armor
$ws.Range("H36").Value = 4943.4287
$ws.Range("I36").Value = 4943.4287
$ws.Range("K36").Value = 4943.4287
$ws.Range("M36").Value = -4597.4287

# Row 97: High Steel Ingot
$ws.Range("H97").Value = 1042.3077
$ws.Range("I97").Value = 968.5454999999999
$ws.Range("K97").Value = 968.5454999999999
$ws.Range("M97").Value = -472.5454999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Iron Ingot
$ws.Range("H20").Value = 2749
$ws.Range("I20").Value = 2749
$ws.Range("K20").Value = 2749
$ws.Range("M20").Value = -2502

# Row 80: Titanium Ingot
$ws.Range("H80").Value = 175
$ws.Range("I80").Value = 260
$ws.Range("J80").Value = 90
$ws.Range("K80").Value = 260
$ws.Range("L80").Value = 90
$ws.Range("M80").Value = 738
$ws.Range("N80").Value = -2086

# Row 83: Titanium Ingot
$ws.Range("H83").Value = 175
$ws.Range("I83").Value = 260
$ws.Range("J83").Value = 90
$ws.Range("K83").Value = 1300
$ws.Range("L83").Value = 450
$ws.Range("M83").Value = 3692
$ws.Range("N83").Value = -10434

# Row 86: Adamantite Nugget
$ws.Range("H86").Value = 1327.2858
$ws.Range("I86").Value = 1308.5
$ws.Range("K86").Value = 1308.5
$ws.Range("M86").Value = -185.5

# Row 89: Adamantite Nugget
$ws.Range("H89").Value = 1327.2858
$ws.Range("I89").Value = 1308.5
$ws.Range("K89").Value = 6542.5
$ws.Range("M89").Value = -926.5

$ws = $wb.Worksheets.Item("CRP")
# Row 99: Pine Lumber
$ws.Range("H99").Value = 2552
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# Row 122: Horse Chestnut Lumber
$ws.Range("H122").Value = 1799.8
$ws.Range("I122").Value = 1799.8
$ws.Range("K122").Value = 5399.4
$ws.Range("M122").Value = -2949.4

# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 2552
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 5195
$ws.Range("I132").Value = 5992.6665
$ws.Range("K132").Value = 17977.9995
$ws.Range("M132").Value = -15447.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 33: Chicken Stock
$ws.Range("H33").Value = 450
$ws.Range("J33").Value = 450
$ws.Range("L33").Value = 2700
$ws.Range("N33").Value = -3266

# Row 44: Knight's Bread
$ws.Range("H44").Value = 398.2143
$ws.Range("I44").Value = 149
$ws.Range("K44").Value = 447
$ws.Range("M44").Value = -49

# Row 98: Rice Vinegar
$ws.Range("H98").Value = 3999.3333
$ws.Range("I98").Value = 3995
$ws.Range("J98").Value = 4001.5
$ws.Range("K98").Value = 11985
$ws.Range("L98").Value = 12004.5
$ws.Range("M98").Value = -10487
$ws.Range("N98").Value = -15000.5

# Row 113: Night Vinegar
$ws.Range("H113").Value = 1456.3334
$ws.Range("I113").Value = 1314.6666
$ws.Range("J113").Value = 1598
$ws.Range("K113").Value = 3943.9998
$ws.Range("L113").Value = 4794
$ws.Range("M113").Value = -1773.9998
$ws.Range("N113").Value = -9134

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Durium Ingot
$ws.Range("H102").Value = 1962.9286
$ws.Range("I102").Value = 1943.1111
$ws.Range("J102").Value = 1998.6
$ws.Range("K102").Value = 1943.1111
$ws.Range("L102").Value = 1998.6
$ws.Range("M102").Value = -321.1111000000001
$ws.Range("N102").Value = -5242.6

# Row 113: Manasilver Nugget
$ws.Range("H113").Value = 999.0909

# Row 134: Ihuykanite
$ws.Range("H134").Value = 39220.832
$ws.Range("J134").Value = 39220.832
$ws.Range("L134").Value = 117662.496
$ws.Range("N134").Value = -122732.496

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 8636.157999999999
$ws.Range("I7").Value = 8333
$ws.Range("J7").Value = 8776.076999999999
$ws.Range("K7").Value = 8333
$ws.Range("L7").Value = 8776.076999999999
$ws.Range("M7").Value = -8221
$ws.Range("N7").Value = -9000.076999999999

# Row 82: Dragon Leather
$ws.Range("H82").Value = 2499.5
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2639

# Row 85: Dragon Leather
$ws.Range("H85").Value = 2499.5
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1752

# Row 126: Saiga Leather
$ws.Range("H126").Value = 8636.157999999999
$ws.Range("I126").Value = 8333
$ws.Range("J126").Value = 8776.076999999999
$ws.Range("K126").Value = 24999
$ws.Range("L126").Value = 26328.231
$ws.Range("M126").Value = -22529
$ws.Range("N126").Value = -31268.231

$ws = $wb.Worksheets.Item("WVR")
# Row 54: Woolen Tights
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Row 81: Crawler Silk
$ws.Range("H81").Value = 835346.0600000001
$ws.Range("I81").Value = 1608.6666
$ws.Range("K81").Value = 3217.3332
$ws.Range("M81").Value = -2156.3332

# Row 84: Crawler Silk
$ws.Range("H84").Value = 835346.0600000001
$ws.Range("I84").Value = 1608.6666
$ws.Range("K84").Value = 16086.666
$ws.Range("M84").Value = -10782.666

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 1970.7142
$ws.Range("I132").Value = 2057.9167
$ws.Range("J132").Value = 1447.5
$ws.Range("K132").Value = 6173.750100000001
$ws.Range("L132").Value = 4342.5
$ws.Range("M132").Value = -3643.750100000001
$ws.Range("N132").Value = -9402.5
